$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 49, pushing the existing rows 49-56 down to 52-59.
$ws.Rows.Item(49).Resize(3).Insert()

# The new rows describe a newer weekly price report (2021-10-05, day value 44474)
# for "Provincia del Elquí" sourced Chirimoya, same shape as the other
# "Provincia del Elquí" rows already in the sheet (Especial/Primera/Segunda).

# Row 49: Especial
$ws.Cells.Item(49, 1).Value = 8
$ws.Cells.Item(49, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(49, 3).Value = "Coquimbo"
$ws.Cells.Item(49, 4).Value = 44474
$ws.Cells.Item(49, 5).Value = 4
$ws.Cells.Item(49, 6).Value = "Fruta"
$ws.Cells.Item(49, 7).Value = 100107
$ws.Cells.Item(49, 8).Value = "Otros"
$ws.Cells.Item(49, 9).Value = 100107002
$ws.Cells.Item(49, 10).Value = "Chirimoya"
$ws.Cells.Item(49, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(49, 12).Value = "Especial"
$ws.Cells.Item(49, 13).Value = 400
$ws.Cells.Item(49, 14).Value = 2100
$ws.Cells.Item(49, 15).Value = 2200
$ws.Cells.Item(49, 16).Value = 2150
$ws.Cells.Item(49, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(49, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(49, 19).Value = 2150
$ws.Cells.Item(49, 20).Value = 1

# Row 50: Primera
$ws.Cells.Item(50, 1).Value = 8
$ws.Cells.Item(50, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(50, 3).Value = "Coquimbo"
$ws.Cells.Item(50, 4).Value = 44474
$ws.Cells.Item(50, 5).Value = 4
$ws.Cells.Item(50, 6).Value = "Fruta"
$ws.Cells.Item(50, 7).Value = 100107
$ws.Cells.Item(50, 8).Value = "Otros"
$ws.Cells.Item(50, 9).Value = 100107002
$ws.Cells.Item(50, 10).Value = "Chirimoya"
$ws.Cells.Item(50, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(50, 12).Value = "Primera"
$ws.Cells.Item(50, 13).Value = 300
$ws.Cells.Item(50, 14).Value = 1700
$ws.Cells.Item(50, 15).Value = 1800
$ws.Cells.Item(50, 16).Value = 1750
$ws.Cells.Item(50, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(50, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(50, 19).Value = 1750
$ws.Cells.Item(50, 20).Value = 1

# Row 51: Segunda
$ws.Cells.Item(51, 1).Value = 8
$ws.Cells.Item(51, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(51, 3).Value = "Coquimbo"
$ws.Cells.Item(51, 4).Value = 44474
$ws.Cells.Item(51, 5).Value = 4
$ws.Cells.Item(51, 6).Value = "Fruta"
$ws.Cells.Item(51, 7).Value = 100107
$ws.Cells.Item(51, 8).Value = "Otros"
$ws.Cells.Item(51, 9).Value = 100107002
$ws.Cells.Item(51, 10).Value = "Chirimoya"
$ws.Cells.Item(51, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(51, 12).Value = "Segunda"
$ws.Cells.Item(51, 13).Value = 200
$ws.Cells.Item(51, 14).Value = 1300
$ws.Cells.Item(51, 15).Value = 1400
$ws.Cells.Item(51, 16).Value = 1350
$ws.Cells.Item(51, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(51, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(51, 19).Value = 1350
$ws.Cells.Item(51, 20).Value = 1

# Ensure the date cells keep the date-formatted style used throughout column D.
$ws.Range("D49:D51").NumberFormat = "YYYY-MM-DD HH:MM:SS"
